# Update input setup for storage China
# - Fill in the "subannual" lvl_temporal value for the charger/discharger
#   (Cha/Dis) rows that previously left column D blank.
# - Widen column D slightly to fit the new content.
# - Re-colour the technology-group legend cells in column W (the fill
#   colours assigned to each 3-row technology block were rotated).
# - Leave the final selection on D16 (last edited cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1. Column D ("lvl_temporal"): fill in "subannual" for the charger
#    ("Cha") / discharger ("Dis") rows that were previously blank.
# ---------------------------------------------------------------
$rowsNeedingSubannual = @(2, 4, 5, 7, 8, 10, 11, 13, 14, 16)
foreach ($r in $rowsNeedingSubannual) {
    $ws.Cells.Item($r, 4).Value = "subannual"
}

# ---------------------------------------------------------------
# 2. Column D is now a touch wider to comfortably fit the text.
# ---------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 11.3

# ---------------------------------------------------------------
# 3. Re-colour the merged legend cells in column W. Each technology
#    block (3 rows, merged) keeps its border/alignment/text but gets
#    reassigned to a different fill colour from the existing palette.
#    Colours below are the exact Interior.Color (OLE BGR) values used
#    elsewhere in the sheet for each palette slot.
# ---------------------------------------------------------------
$paletteColor = @{
    4 = 16247773   # theme accent1, lighter 80%
    5 = 13553360   # theme background2, darker 10%
    6 = 6740479    # theme accent4, lighter 40%
    7 = 9359529    # theme accent6, lighter 40%
    8 = 15123099   # theme accent1, lighter 40%
}

$groupFill = @{
    "W2:W4"   = 8
    "W5:W7"   = 7
    "W8:W10"  = 7
    "W11:W13" = 7
    "W14:W16" = 7
    "W17:W19" = 4
    "W20:W22" = 7
    "W23:W25" = 5
    "W26:W28" = 6
}

foreach ($rangeAddr in $groupFill.Keys) {
    $fillSlot = $groupFill[$rangeAddr]
    $ws.Range($rangeAddr).Interior.Color = $paletteColor[$fillSlot]
}

# ---------------------------------------------------------------
# 4. Leave the cursor on the last touched cell, like the source edit.
# ---------------------------------------------------------------
$ws.Range("D16").Select() | Out-Null
